$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# James -> Duncan for the roster entry in row 8 (Keller's first name)
$ws.Range("B8").Value = "Duncan"

# Update the active selection to B9 as recorded in the saved view state
$null = $ws.Range("B9").Select()
